$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 27780908
$ws.Range("I6").Value = 71428760
$ws.Range("J6").Value = 4998.909
$ws.Range("K6").Value = 214286280
$ws.Range("L6").Value = 14996.727
$ws.Range("M6").Value = -214286168
$ws.Range("N6").Value = -15220.727

$ws.Range("H9").Value = 714379.7
$ws.Range("I9").Value = 1666735.6
$ws.Range("J9").Value = 112.75
$ws.Range("K9").Value = 1666735.6
$ws.Range("L9").Value = 112.75
$ws.Range("M9").Value = -1666566.6
$ws.Range("N9").Value = -450.75

$ws.Range("H12").Value = 412.14285
$ws.Range("I12").Value = 150
$ws.Range("K12").Value = 150
$ws.Range("M12").Value = 20

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").Value = 0

$ws.Range("H29").Value = 62500176
$ws.Range("I29").Value = 62500176
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 187500528
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -187500247

$ws.Range("H33").Value = 631.2174
$ws.Range("I33").Value = 166.64285
$ws.Range("K33").Value = 166.64285
$ws.Range("M33").Value = 62.35714999999999

$ws.Range("H38").Value = 3026.5264
$ws.Range("I38").Value = 682.63635
$ws.Range("J38").Value = 6249.375
$ws.Range("K38").Value = 2047.90905
$ws.Range("L38").Value = 18748.125
$ws.Range("M38").Value = -1675.90905
$ws.Range("N38").Value = -19492.125

$ws.Range("H58").Value = 35714950
$ws.Range("I58").Value = 35714950
$ws.Range("K58").Value = 107144850
$ws.Range("M58").Value = -107144700

$ws.Range("H76").Value = 3833.3333
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685

$ws.Range("H79").Value = 3833.3333
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908

$ws.Range("H80").Value = 924.4211
$ws.Range("I80").Value = 244.71428
$ws.Range("K80").Value = 734.14284
$ws.Range("M80").Value = 263.85716

$ws.Range("H83").Value = 924.4211
$ws.Range("I83").Value = 244.71428
$ws.Range("K83").Value = 2202.42852
$ws.Range("M83").Value = 2789.57148

$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = 0

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 0

$ws.Range("H141").Value = 7634.773
$ws.Range("I141").Value = 4716
$ws.Range("K141").Value = 14148
$ws.Range("M141").Value = -8968


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1995.1666
$ws.Range("I45").Value = 2043
$ws.Range("K45").Value = 2043
$ws.Range("M45").Value = -1666

$ws.Range("H74").Value = 3506.875
$ws.Range("I74").Value = 3506.875
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3506.875
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -2632.875

$ws.Range("H77").Value = 3506.875
$ws.Range("I77").Value = 3506.875
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 17534.375
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -13166.375

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 360.17242
$ws.Range("J80").Value = 367.85
$ws.Range("L80").Value = 367.85
$ws.Range("N80").Value = -2363.85

$ws.Range("H83").Value = 360.17242
$ws.Range("J83").Value = 367.85
$ws.Range("L83").Value = 1839.25
$ws.Range("N83").Value = -11823.25

$ws.Range("H92").Value = 29749.334
$ws.Range("J92").Value = 29749.334
$ws.Range("L92").Value = 29749.334
$ws.Range("N92").Value = -34741.334


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 67500
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998

$ws.Range("H71").Value = 67500
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 550
$ws.Range("J68").Value = 550
$ws.Range("L68").Value = 1650
$ws.Range("N68").Value = -3272

$ws.Range("H71").Value = 550
$ws.Range("J71").Value = 550
$ws.Range("L71").Value = 4950
$ws.Range("N71").Value = -13062


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 630.8570999999999
$ws.Range("I2").Value = 66.111115
$ws.Range("K2").Value = 66.111115
$ws.Range("M2").Value = 46.888885

$ws.Range("H63").Value = 29999
$ws.Range("J63").Value = 29999
$ws.Range("L63").Value = 29999
$ws.Range("N63").Value = -31371

$ws.Range("H66").Value = 29999
$ws.Range("J66").Value = 29999
$ws.Range("L66").Value = 89997
$ws.Range("N66").Value = -96861

$ws.Range("H102").Value = 2108
$ws.Range("I102").Value = 2201.8635
$ws.Range("J102").Value = 1849.875
$ws.Range("K102").Value = 2201.8635
$ws.Range("L102").Value = 1849.875
$ws.Range("M102").Value = -579.8634999999999
$ws.Range("N102").Value = -5093.875

$ws.Range("H113").Value = 8888.235000000001
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 10700
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 10700
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -15040


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3549.7837
$ws.Range("I61").Value = 2789.2593
$ws.Range("K61").Value = 2789.2593
$ws.Range("M61").Value = -2587.2593

$ws.Range("H113").Value = 3549.7837
$ws.Range("I113").Value = 2789.2593
$ws.Range("K113").Value = 2789.2593
$ws.Range("M113").Value = -619.2593000000002


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 21598
$ws.Range("I24").Value = 21598
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 21598
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -21368

$ws.Range("H99").Value = 33999
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 2924.375
$ws.Range("I107").Value = 4166.6665
$ws.Range("J107").Value = 2179
$ws.Range("K107").Value = 12499.9995
$ws.Range("L107").Value = 6537
$ws.Range("M107").Value = -10579.9995
$ws.Range("N107").Value = -10377

